$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62 (previously the last row) loses the "last row" date format and
# reverts to the regular date/time format used by every other row above it.
$ws.Range("A62").NumberFormat = $ws.Range("A61").NumberFormat

# Row 63 is the new daily entry; it inherits the "last row" date format
# that row 62 used to have.
$ws.Range("A63").NumberFormat = "YYYY-MM-DD"
$ws.Range("A63").Value = 45803
$ws.Range("B63").Value = 264
$ws.Range("C63").Value = 273
$ws.Range("D63").Value = 268
